$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Progress update as of 04-Nov-2025: "PERIOD TO EXPIRE" (H) drops by one day
# and "LAST UPDATE" (I) moves from 03-Nov-2025 to 04-Nov-2025 for rows 3-7.
#
# I column cells hold literal text dates (not real Excel dates), so a plain
# .Value assignment of "04-Nov-2025" would be auto-parsed into a date serial
# and force a new number-formatted style. To keep the cells as plain text
# with their original style, write the text via a formula and then convert
# that formula to a static value with Copy / PasteSpecial (xlPasteValues).

$ws.Range("H3").Value = 64
$ws.Range("I3").Formula = "=""04-Nov-2025"""
$ws.Range("I3").Copy()
$ws.Range("I3").PasteSpecial(-4163)

$ws.Range("H4").Value = 631
$ws.Range("I4").Formula = "=""04-Nov-2025"""
$ws.Range("I4").Copy()
$ws.Range("I4").PasteSpecial(-4163)

$ws.Range("H5").Value = 630
$ws.Range("I5").Formula = "=""04-Nov-2025"""
$ws.Range("I5").Copy()
$ws.Range("I5").PasteSpecial(-4163)

$ws.Range("H6").Value = 631
$ws.Range("I6").Formula = "=""04-Nov-2025"""
$ws.Range("I6").Copy()
$ws.Range("I6").PasteSpecial(-4163)

$ws.Range("H7").Value = 630
$ws.Range("I7").Formula = "=""04-Nov-2025"""
$ws.Range("I7").Copy()
$ws.Range("I7").PasteSpecial(-4163)

$excel.CutCopyMode = $false
